$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.584298133850098
$ws.Range("B1").Value = 2.301728248596191
$ws.Range("C1").Value = 2.702759981155396
$ws.Range("D1").Value = 3.266005992889404
$ws.Range("E1").Value = 1.648396253585815
